# Scheduled refresh of Asura Profits leve pricing data.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N)
# for the affected leve rows across the ALC, ARM, BSM, CRP, CUL, GSM,
# LTW and WVR sheets with freshly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1797.5922
$ws.Range("I15").Value = 1797.5922
$ws.Range("K15").Value = 5392.7766
$ws.Range("M15").Value = -5223.7766

$ws.Range("H76").Value = 4509.364
$ws.Range("I76").Value = 4829
$ws.Range("J76").Value = 3950
$ws.Range("K76").Value = 4829
$ws.Range("L76").Value = 3950
$ws.Range("M76").Value = -4514
$ws.Range("N76").Value = -4580

$ws.Range("H79").Value = 4509.364
$ws.Range("I79").Value = 4829
$ws.Range("J79").Value = 3950
$ws.Range("K79").Value = 4829
$ws.Range("L79").Value = 3950
$ws.Range("M79").Value = -3737
$ws.Range("N79").Value = -6134

$ws.Range("H80").Value = 7721.1875
$ws.Range("I80").Value = 1104.5
$ws.Range("J80").Value = 14337.875
$ws.Range("K80").Value = 3313.5
$ws.Range("L80").Value = 43013.625
$ws.Range("M80").Value = -2315.5
$ws.Range("N80").Value = -45009.625

$ws.Range("H83").Value = 7721.1875
$ws.Range("I83").Value = 1104.5
$ws.Range("J83").Value = 14337.875
$ws.Range("K83").Value = 9940.5
$ws.Range("L83").Value = 129040.875
$ws.Range("M83").Value = -4948.5
$ws.Range("N83").Value = -139024.875

$ws.Range("H132").Value = 2621.111
$ws.Range("I132").Value = 2257.2068
$ws.Range("J132").Value = 3280.6875
$ws.Range("K132").Value = 6771.6204
$ws.Range("L132").Value = 9842.0625
$ws.Range("M132").Value = -4241.6204
$ws.Range("N132").Value = -14902.0625

$ws.Range("H137").Value = 1162.0156
$ws.Range("I137").Value = 953.66
$ws.Range("J137").Value = 1906.1428
$ws.Range("K137").Value = 2860.98
$ws.Range("L137").Value = 5718.428400000001
$ws.Range("M137").Value = -310.98
$ws.Range("N137").Value = -10818.4284

$ws.Range("H138").Value = 2636.878
$ws.Range("I138").Value = 1672.7885
$ws.Range("J138").Value = 4307.967
$ws.Range("K138").Value = 5018.3655
$ws.Range("L138").Value = 12923.901
$ws.Range("M138").Value = 121.6345000000001
$ws.Range("N138").Value = -23203.901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1144.7561
$ws.Range("I61").Value = 1120.3823
$ws.Range("K61").Value = 1120.3823
$ws.Range("M61").Value = -908.3823

$ws.Range("H74").Value = 1073.1351
$ws.Range("I74").Value = 1090.9412
$ws.Range("J74").Value = 871.3333
$ws.Range("K74").Value = 1090.9412
$ws.Range("L74").Value = 871.3333
$ws.Range("M74").Value = -216.9412
$ws.Range("N74").Value = -2619.3333

$ws.Range("H77").Value = 1073.1351
$ws.Range("I77").Value = 1090.9412
$ws.Range("J77").Value = 871.3333
$ws.Range("K77").Value = 5454.706
$ws.Range("L77").Value = 4356.6665
$ws.Range("M77").Value = -1086.706
$ws.Range("N77").Value = -13092.6665

$ws.Range("H123").Value = 24314.715
$ws.Range("J123").Value = 24314.715
$ws.Range("L123").Value = 24314.715
$ws.Range("N123").Value = -34114.715

$ws.Range("H132").Value = 2807.8333
$ws.Range("I132").Value = 2839.1428
$ws.Range("J132").Value = 2794.9412
$ws.Range("K132").Value = 8517.428400000001
$ws.Range("L132").Value = 8384.8236
$ws.Range("M132").Value = -5987.428400000001
$ws.Range("N132").Value = -13444.8236

$ws.Range("H136").Value = 1144.7561
$ws.Range("I136").Value = 1120.3823
$ws.Range("K136").Value = 3361.1469
$ws.Range("M136").Value = -811.1468999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3037.8572
$ws.Range("I134").Value = 3196.818
$ws.Range("J134").Value = 2935
$ws.Range("K134").Value = 9590.454000000002
$ws.Range("L134").Value = 8805
$ws.Range("M134").Value = -7055.454000000002
$ws.Range("N134").Value = -13875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2605.12
$ws.Range("I31").Value = 1754.1765
$ws.Range("J31").Value = 4413.375
$ws.Range("K31").Value = 1754.1765
$ws.Range("L31").Value = 4413.375
$ws.Range("M31").Value = -1459.1765
$ws.Range("N31").Value = -5003.375

$ws.Range("H34").Value = 2605.12
$ws.Range("I34").Value = 1754.1765
$ws.Range("J34").Value = 4413.375
$ws.Range("K34").Value = 1754.1765
$ws.Range("L34").Value = 4413.375
$ws.Range("M34").Value = -1552.1765
$ws.Range("N34").Value = -4817.375

$ws.Range("H60").Value = 38899.668
$ws.Range("J60").Value = 38899.668
$ws.Range("L60").Value = 38899.668
$ws.Range("N60").Value = -39921.668

$ws.Range("H99").Value = 5700
$ws.Range("I99").Value = 6000
$ws.Range("K99").Value = 6000
$ws.Range("M99").Value = -4502

$ws.Range("H122").Value = 2561.4688
$ws.Range("I122").Value = 2688.5833
$ws.Range("J122").Value = 2180.125
$ws.Range("K122").Value = 8065.749899999999
$ws.Range("L122").Value = 6540.375
$ws.Range("M122").Value = -5615.749899999999
$ws.Range("N122").Value = -11440.375

$ws.Range("H126").Value = 5700
$ws.Range("I126").Value = 6000
$ws.Range("K126").Value = 18000
$ws.Range("M126").Value = -15530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1600
$ws.Range("I59").Value = 1600
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 4800
$ws.Range("L59").ClearContents()
$ws.Range("M59").Value = -4260
$ws.Range("N59").Value = 0

$ws.Range("H132").Value = 2331.111
$ws.Range("I132").Value = 1480
$ws.Range("J132").Value = 3395
$ws.Range("K132").Value = 13320
$ws.Range("L132").Value = 30555
$ws.Range("M132").Value = -10790
$ws.Range("N132").Value = -35615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = 0

$ws.Range("H70").Value = 7030.45
$ws.Range("I70").Value = 6466.6665
$ws.Range("K70").Value = 6466.6665
$ws.Range("M70").Value = -6196.6665

$ws.Range("H73").Value = 7030.45
$ws.Range("I73").Value = 6466.6665
$ws.Range("K73").Value = 6466.6665
$ws.Range("M73").Value = -5530.6665

$ws.Range("H80").Value = 2630.5
$ws.Range("I80").Value = 2672.1428
$ws.Range("J80").Value = 2533.3333
$ws.Range("K80").Value = 2672.1428
$ws.Range("L80").Value = 2533.3333
$ws.Range("M80").Value = -1674.1428
$ws.Range("N80").Value = -4529.3333

$ws.Range("H83").Value = 2630.5
$ws.Range("I83").Value = 2672.1428
$ws.Range("J83").Value = 2533.3333
$ws.Range("K83").Value = 13360.714
$ws.Range("L83").Value = 12666.6665
$ws.Range("M83").Value = -8368.714
$ws.Range("N83").Value = -22650.6665

$ws.Range("H102").Value = 5437.5625
$ws.Range("I102").Value = 4500.1665
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 4500.1665
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -2878.1665
$ws.Range("N102").Value = -9244

$ws.Range("H132").Value = 2278.7036
$ws.Range("I132").Value = 1487.7778
$ws.Range("J132").Value = 3860.5557
$ws.Range("K132").Value = 4463.3334
$ws.Range("L132").Value = 11581.6671
$ws.Range("M132").Value = -1933.3334
$ws.Range("N132").Value = -16641.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1723.1333
$ws.Range("I82").Value = 1554.9
$ws.Range("J82").Value = 2059.6
$ws.Range("K82").Value = 1554.9
$ws.Range("L82").Value = 2059.6
$ws.Range("M82").Value = -1193.9
$ws.Range("N82").Value = -2781.6

$ws.Range("H85").Value = 1723.1333
$ws.Range("I85").Value = 1554.9
$ws.Range("J85").Value = 2059.6
$ws.Range("K85").Value = 1554.9
$ws.Range("L85").Value = 2059.6
$ws.Range("M85").Value = -306.9000000000001
$ws.Range("N85").Value = -4555.6

$ws.Range("H94").Value = 20165
$ws.Range("J94").Value = 20165
$ws.Range("L94").Value = 20165
$ws.Range("N94").Value = -21517

$ws.Range("H122").Value = 12861236
$ws.Range("I122").Value = 10420577
$ws.Range("J122").Value = 18186310
$ws.Range("K122").Value = 31261731
$ws.Range("L122").Value = 54558930
$ws.Range("M122").Value = -31259281
$ws.Range("N122").Value = -54563830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4500.5
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4858
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4858
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6106

$ws.Range("H65").Value = 4500.5
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4858
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 24290
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -30530
